# enemyDatabase.xlsx — Critical / crit-chance spell property pass.
#
# Renames a few "style mod" / spell-name strings on the single
# "enemyDatabase" sheet:
#   - Slime  (row 2) "Spell x (root)"  : "poke"  -> "weakpnt"
#   - Tanuki (row 3) "Spell1 (root)"   : "whiff" -> "spear"
#   - Tanuki (row 3) "Spell x (style)" : "null"  -> "aimed"
#
# Order matters here only insofar as it determines the order brand-new
# shared-string entries are appended in; it otherwise has no visible
# effect, so we just go row by row, top to bottom, left to right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tanuki's root spell ("Spell1 (root)", column O) was "whiff", now "spear".
$ws.Range("O3").Value = "spear"

# Tanuki's style mod ("Spell x (style)", column T) was "null", now "aimed".
$ws.Range("T3").Value = "aimed"

# Slime's secondary spell root ("Spell x (root)", column R) was "poke",
# now "weakpnt".
$ws.Range("R2").Value = "weakpnt"

# The author's cursor ended up on R2 when the workbook was saved.
$ws.Range("R2").Select()
